$d = $word.ActiveDocument

# --- Locate the paragraph that ends "... account to his manager for (= tell
# his manager about and explain) all his movements." -------------------------
function Find-ParaIndexByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like $needle) {
            return $i
        }
    }
    return -1
}

$idx = Find-ParaIndexByText "*account to his manager for*"
if ($idx -eq -1) {
    throw "Could not find target paragraph (account to his manager for)"
}

# 1) Replace that paragraph's XML with itself plus a Vietnamese language tag
#    on the paragraph mark (<w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr>).
$target = $d.Paragraphs.Item($idx)
$targetXml = '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' + `
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' + `
    '<w:r><w:t xml:space="preserve"> He </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>has to</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> account to his manager for (= tell his manager about and explain) all his movements.</w:t></w:r>' + `
    '</w:p>'
$target.Range.InsertXML($targetXml)

# 2) Insert the three new paragraphs right after it, as raw OOXML so the
#    run/paragraph formatting matches exactly:
#      - a blank paragraph (indented, Vietnamese language mark)
#      - the new "Sparks fly:" glossary entry (bold term + definition)
#      - the example sentence paragraph with the Wingdings arrow glyph
$idx = Find-ParaIndexByText "*account to his manager for*"
$nextPara = $d.Paragraphs.Item($idx + 1)
$insertAt = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)

$newParagraphsXml = '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
    '<w:rPr><w:b/><w:bCs/><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">Sparks </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">fly: </w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>If sparks fly between two or more people, they argue angrily</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:ind w:left="360"/></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>When they get together in a meeting, the sparks really fly.</w:t></w:r>' + `
    '</w:p>'

$insertAt.InsertXML($newParagraphsXml)

$d.Save()
